# Daily attendance processing - 2026-01-09 11:55:40
# Reorders the "Recorded By" (column G) names so that the trailing
# comma-separated entry is moved to the front of the list, for any
# cell whose value begins with "System, " or "system, System, ".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7
    $val = $cell.Value2

    if ($null -eq $val) { continue }

    $text = [string]$val

    if ($text.StartsWith("System, ") -or $text.StartsWith("system, System, ")) {
        $parts = $text.Split(",")
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }
        $last = $parts[$parts.Length - 1]
        $rest = $parts[0..($parts.Length - 2)]
        $newParts = @($last) + $rest
        $newText = [string]::Join(", ", $newParts)
        $cell.Value = $newText
    }
}
